$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set B/C (text) columns first where coin identity changed
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('B27').Value = 'BitcoinCash'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'

# Set D (price) columns as text to avoid numeric auto-conversion
$d = $ws.Range('D2')
$d.NumberFormat = '@'
$d.Value = '26.335.68'
$d.Style = 'Normal'
$d = $ws.Range('D3')
$d.NumberFormat = '@'
$d.Value = '1.717.41'
$d.Style = 'Normal'
$d = $ws.Range('D4')
$d.NumberFormat = '@'
$d.Value = '0.9998'
$d.Style = 'Normal'
$d = $ws.Range('D5')
$d.NumberFormat = '@'
$d.Value = '239.23'
$d.Style = 'Normal'
$d = $ws.Range('D7')
$d.NumberFormat = '@'
$d.Value = '0.4736'
$d.Style = 'Normal'
$d = $ws.Range('D8')
$d.NumberFormat = '@'
$d.Value = '0.2629'
$d.Style = 'Normal'
$d = $ws.Range('D9')
$d.NumberFormat = '@'
$d.Value = '0.06208'
$d.Style = 'Normal'
$d = $ws.Range('D10')
$d.NumberFormat = '@'
$d.Value = '1.715.01'
$d.Style = 'Normal'
$d = $ws.Range('D11')
$d.NumberFormat = '@'
$d.Value = '0.07057'
$d.Style = 'Normal'
$d = $ws.Range('D12')
$d.NumberFormat = '@'
$d.Value = '15.32'
$d.Style = 'Normal'
$d = $ws.Range('D13')
$d.NumberFormat = '@'
$d.Value = '4.415'
$d.Style = 'Normal'
$d = $ws.Range('D14')
$d.NumberFormat = '@'
$d.Value = '0.5891'
$d.Style = 'Normal'
$d = $ws.Range('D15')
$d.NumberFormat = '@'
$d.Value = '76.00'
$d.Style = 'Normal'
$d = $ws.Range('D17')
$d.NumberFormat = '@'
$d.Value = '1.001'
$d.Style = 'Normal'
$d = $ws.Range('D18')
$d.NumberFormat = '@'
$d.Value = '26.329.91'
$d.Style = 'Normal'
$d = $ws.Range('D19')
$d.NumberFormat = '@'
$d.Value = '0.000006818'
$d.Style = 'Normal'
$d = $ws.Range('D20')
$d.NumberFormat = '@'
$d.Value = '11.53'
$d.Style = 'Normal'
$d = $ws.Range('D21')
$d.NumberFormat = '@'
$d.Value = '1.936.04'
$d.Style = 'Normal'
$d = $ws.Range('D22')
$d.NumberFormat = '@'
$d.Value = '4.536'
$d.Style = 'Normal'
$d = $ws.Range('D23')
$d.NumberFormat = '@'
$d.Value = '8.747'
$d.Style = 'Normal'
$d = $ws.Range('D24')
$d.NumberFormat = '@'
$d.Value = '5.314'
$d.Style = 'Normal'
$d = $ws.Range('D25')
$d.NumberFormat = '@'
$d.Value = '134.34'
$d.Style = 'Normal'
$d = $ws.Range('D26')
$d.NumberFormat = '@'
$d.Value = '15.26'
$d.Style = 'Normal'
$d = $ws.Range('D27')
$d.NumberFormat = '@'
$d.Value = '108.27'
$d.Style = 'Normal'
$d = $ws.Range('D28')
$d.NumberFormat = '@'
$d.Value = '1.400'
$d.Style = 'Normal'
$d = $ws.Range('D29')
$d.NumberFormat = '@'
$d.Value = '1.749'
$d.Style = 'Normal'
$d = $ws.Range('D30')
$d.NumberFormat = '@'
$d.Value = '3.992'
$d.Style = 'Normal'
$d = $ws.Range('D31')
$d.NumberFormat = '@'
$d.Value = '3.681'
$d.Style = 'Normal'
$d = $ws.Range('D32')
$d.NumberFormat = '@'
$d.Value = '0.07753'
$d.Style = 'Normal'
$d = $ws.Range('D33')
$d.NumberFormat = '@'
$d.Value = '0.04429'
$d.Style = 'Normal'
$d = $ws.Range('D34')
$d.NumberFormat = '@'
$d.Value = '2.613'
$d.Style = 'Normal'
$d = $ws.Range('D35')
$d.NumberFormat = '@'
$d.Value = '0.9757'
$d.Style = 'Normal'
$d = $ws.Range('D36')
$d.NumberFormat = '@'
$d.Value = '0.6179'
$d.Style = 'Normal'
$d = $ws.Range('D37')
$d.NumberFormat = '@'
$d.Value = '0.9329'
$d.Style = 'Normal'
$d = $ws.Range('D38')
$d.NumberFormat = '@'
$d.Value = '112.87'
$d.Style = 'Normal'
$d = $ws.Range('D39')
$d.NumberFormat = '@'
$d.Value = '2.421'
$d.Style = 'Normal'
$d = $ws.Range('D40')
$d.NumberFormat = '@'
$d.Value = '1.911'
$d.Style = 'Normal'
$d = $ws.Range('D41')
$d.NumberFormat = '@'
$d.Value = '1.000'
$d.Style = 'Normal'
$d = $ws.Range('D42')
$d.NumberFormat = '@'
$d.Value = '0.01473'
$d.Style = 'Normal'
$d = $ws.Range('D43')
$d.NumberFormat = '@'
$d.Value = '5.326'
$d.Style = 'Normal'
$d = $ws.Range('D44')
$d.NumberFormat = '@'
$d.Value = '0.3801'
$d.Style = 'Normal'
$d = $ws.Range('D45')
$d.NumberFormat = '@'
$d.Value = '0.1166'
$d.Style = 'Normal'
$d = $ws.Range('D46')
$d.NumberFormat = '@'
$d.Value = '6.296'
$d.Style = 'Normal'
$d = $ws.Range('D47')
$d.NumberFormat = '@'
$d.Value = '0.05282'
$d.Style = 'Normal'
$d = $ws.Range('D49')
$d.NumberFormat = '@'
$d.Value = '7.690'
$d.Style = 'Normal'
$d = $ws.Range('D50')
$d.NumberFormat = '@'
$d.Value = '50.65'
$d.Style = 'Normal'
$d = $ws.Range('D51')
$d.NumberFormat = '@'
$d.Value = '0.3357'
$d.Style = 'Normal'

# Set E (volume %) columns
$ws.Range('E2').Value = '  +2.98%  '
$ws.Range('E3').Value = '  +3.20%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('E5').Value = '  +1.41%  '
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('E7').Value = '  -1.22%  '
$ws.Range('E8').Value = '  +0.49%  '
$ws.Range('E9').Value = '  +0.97%  '
$ws.Range('E10').Value = '  +3.09%  '
$ws.Range('E11').Value = '  -0.34%  '
$ws.Range('E12').Value = '  +3.91%  '
$ws.Range('E13').Value = '  +0.98%  '
$ws.Range('E14').Value = '  -0.39%  '
$ws.Range('E15').Value = '  +2.13%  '
$ws.Range('E16').Value = '  +0.01%  '
$ws.Range('E17').Value = '  +0.11%  '
$ws.Range('E18').Value = '  +3.03%  '
$ws.Range('E19').Value = '  +0.94%  '
$ws.Range('E20').Value = '  +1.22%  '
$ws.Range('E21').Value = '  +3.02%  '
$ws.Range('E22').Value = '  +2.44%  '
$ws.Range('E23').Value = '  +1.08%  '
$ws.Range('E24').Value = '  +0.45%  '
$ws.Range('E25').Value = '  -0.06%  '
$ws.Range('E26').Value = '  +1.39%  '
$ws.Range('E27').Value = '  +3.52%  '
$ws.Range('E28').Value = '  -0.09%  '
$ws.Range('E29').Value = '  +4.00%  '
$ws.Range('E30').Value = '  +0.86%  '
$ws.Range('E31').Value = '  +0.78%  '
$ws.Range('E32').Value = '  +1.77%  '
$ws.Range('E33').Value = '  +2.62%  '
$ws.Range('E34').Value = '  -0.25%  '
$ws.Range('E35').Value = '  +2.88%  '
$ws.Range('E36').Value = '  +1.11%  '
$ws.Range('E37').Value = '  +9.79%  '
$ws.Range('E38').Value = '  +15.21%  '
$ws.Range('E39').Value = '  -7.24%  '
$ws.Range('E40').Value = '  +2.39%  '
$ws.Range('E41').Value = '  -0.01%  '
$ws.Range('E42').Value = '  -1.62%  '
$ws.Range('E43').Value = '  +13.67%  '
$ws.Range('E44').Value = '  +1.12%  '
$ws.Range('E45').Value = '  +4.21%  '
$ws.Range('E46').Value = '  +1.37%  '
$ws.Range('E47').Value = '  +0.42%  '
$ws.Range('E48').Value = '  +2.89%  '
$ws.Range('E49').Value = '  +5.45%  '
$ws.Range('E50').Value = '  +1.00%  '
$ws.Range('E51').Value = '  +0.76%  '
